# Updated capital structure database
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "D"  = -0.0119
    "E"  = -0.008529999999999999
    "G"  = 0.1582213029989659
    "H"  = 0.1582213029989659
    "I"  = 0.1737331954498449
    "J"  = 0.145996843193817
    "K"  = 1.94
    "L"  = 0.2006204756980352
    "M"  = 1.875
    "N"  = 0.04040948275862069
    "O"  = 0.9664948453608248
    "P"  = 1.54
    "Q"  = 0.0331896551724138
    "R"  = 0.7938144329896908
    "S"  = 0.335
    "T"  = 0.1786666666666666
    "U"  = 3.39
    "V"  = 0.07306034482758621
    "W"  = 0.07822580645161289
    "X"  = 0.02074394029738349
    "Y"  = 0.05748186615422941
    "Z"  = 0.7325757575757575
    "AA" = 0.1069537480063796
    "AB" = 0.0209779607983485
    "AC" = 0.08597578720803108
    "AD" = 0.945
    "AE" = 0
    "AF" = 0.945
    "AG" = -2.445
    "AH" = 0.01995986904636181
    "AI" = 0.02550951545417735
    "AJ" = -0.05562507109543852
    "AK" = -0.0726489377507057
    "AL" = 0.089
    "AM" = -0.652
    "AN" = 0.3987341772151898
    "AO" = 18.87640449438202
    "AP" = -1.031645569620253
    "AQ" = -2.576687116564417
}

foreach ($row in @(2, 3)) {
    foreach ($col in $newValues.Keys) {
        $ws.Range("$col$row").Value = $newValues[$col]
    }
}
